$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to make room for two new reporting quarters
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from the (now shifted) data columns F:M back onto the new D:E columns
# so the new columns inherit the correct date / number styles used throughout the sheet
$ws.Range("F5:M102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D = period ending 2018-06-30, E = period ending 2017-12-31)
$ws.Range("D7").Value = 43281
$ws.Range("E7").Value = 43100
$ws.Range("D8").Value = 14700
$ws.Range("E8").Value = 13100
$ws.Range("D9").Value = 18500
$ws.Range("E9").Value = 27800
$ws.Range("D10").Value = -3800
$ws.Range("E10").Value = -14700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = 67300
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 36000
$ws.Range("E17").Value = 124400
$ws.Range("D18").Value = -21300
$ws.Range("E18").Value = -111300
$ws.Range("D20").Value = 1000
$ws.Range("E20").Value = 1900
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -20400
$ws.Range("E23").Value = -109400
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -20400
$ws.Range("E26").Value = -109600
$ws.Range("D27").Value = -18500
$ws.Range("E27").Value = -104100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1000
$ws.Range("E32").Value = -1900
$ws.Range("D33").Value = -18500
$ws.Range("E33").Value = -104100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -18500
$ws.Range("E35").Value = -104100
$ws.Range("D38").Value = 43281
$ws.Range("E38").Value = 43100
$ws.Range("D41").Value = 4400
$ws.Range("E41").Value = 15400
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("D43").Value = 13200
$ws.Range("E43").Value = 13200
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 62400
$ws.Range("E45").Value = 70000
$ws.Range("D46").Value = 79900
$ws.Range("E46").Value = 98600
$ws.Range("D47").Value = 100300
$ws.Range("E47").Value = 102400
$ws.Range("D48").Value = 14800
$ws.Range("E48").Value = 15400
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8100
$ws.Range("E52").Value = 8500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 203200
$ws.Range("E54").Value = 225000
$ws.Range("D57").Value = 46200
$ws.Range("E57").Value = 48500
$ws.Range("D58").Value = 6200
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 25800
$ws.Range("E59").Value = 27300
$ws.Range("D60").Value = 78200
$ws.Range("E60").Value = 75800
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 24700
$ws.Range("E62").Value = 25500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 76600
$ws.Range("E66").Value = 77400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -190900
$ws.Range("E72").Value = -172300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 126600
$ws.Range("E76").Value = 147600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43281
$ws.Range("E80").Value = 43100
$ws.Range("D81").Value = -18500
$ws.Range("E81").Value = -104100
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("E89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 0

# Fix a handful of cells where the quarter realignment also changed availability of data (NA vs numeric)
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = "NA"
$ws.Range("L89").Value = "NA"
$ws.Range("M89").Value = -1800
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = "NA"
$ws.Range("L94").Value = "NA"
$ws.Range("M94").Value = -6200
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = "NA"
$ws.Range("L100").Value = "NA"
$ws.Range("M100").Value = 16800
